$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = "Unknown"
$ws.Range("H2").Value = "2025-05-12 06:15"

$ws.Range("G3").Value = "UNKNOWN"
$ws.Range("H3").Value = "2025-05-12 06:15"
